$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row at row 113 (shifting the existing rows 113:174 down to 114:175).
$ws.Rows.Item(113).Insert()

# Populate the newly inserted row 113 with the new record's data.
$ws.Cells.Item(113, 1).Value = 11
$ws.Cells.Item(113, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(113, 3).Value = 'Bíobío'
$ws.Cells.Item(113, 4).Value = 44719
$ws.Cells.Item(113, 5).Value = 8
$ws.Cells.Item(113, 6).Value = 'Fruta'
$ws.Cells.Item(113, 7).Value = 100108
$ws.Cells.Item(113, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(113, 9).Value = 100108005
$ws.Cells.Item(113, 10).Value = 'Piña'
$ws.Cells.Item(113, 11).Value = 'Sin especificar'
$ws.Cells.Item(113, 12).Value = 'Segunda'
$ws.Cells.Item(113, 13).Value = 140
$ws.Cells.Item(113, 14).Value = 15000
$ws.Cells.Item(113, 15).Value = 16000
$ws.Cells.Item(113, 16).Value = 15571
$ws.Cells.Item(113, 17).Value = '$/caja 14 unidades'
$ws.Cells.Item(113, 18).Value = 'Ecuador'
$ws.Cells.Item(113, 19).Value = 1112
$ws.Cells.Item(113, 20).Value = 14
